# Updates the cryptos price list with freshly scraped Price / Volume(1h)
# values (GitHub Actions scrape refresh), and fixes the Mantle / RocketPoolETH
# row ordering in rows 47-48.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'29.202.09"
$ws.Range("E2").Value = "  -0.05%  "
$ws.Range("D3").Value = "'1.848.60"
$ws.Range("E3").Value = "  -0.64%  "
$ws.Range("D4").Value = "'0.9988"
$ws.Range("E4").Value = "  -0.10%  "
$ws.Range("D5").Value = "'246.12"
$ws.Range("E5").Value = "  +1.85%  "
$ws.Range("D6").Value = "'0.6984"
$ws.Range("E6").Value = "  -0.65%  "
$ws.Range("D7").Value = "'0.9995"
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("D8").Value = "'0.07722"
$ws.Range("E8").Value = "  -1.19%  "
$ws.Range("E9").Value = "  -1.64%  "
$ws.Range("D10").Value = "'23.57"
$ws.Range("E10").Value = "  -1.21%  "
$ws.Range("D11").Value = "'0.07822"
$ws.Range("E11").Value = "  +0.25%  "
$ws.Range("D12").Value = "'93.44"
$ws.Range("E12").Value = "  +0.90%  "
$ws.Range("D13").Value = "'5.131"
$ws.Range("E13").Value = "  +0.12%  "
$ws.Range("D14").Value = "'1.841.08"
$ws.Range("E14").Value = "  -0.97%  "
$ws.Range("D15").Value = "'0.6865"
$ws.Range("E15").Value = "  -0.45%  "
$ws.Range("D16").Value = "'6.655"
$ws.Range("E16").Value = "  +1.63%  "
$ws.Range("D17").Value = "'0.000008316"
$ws.Range("E17").Value = "  -1.46%  "
$ws.Range("D18").Value = "'29.185.81"
$ws.Range("E18").Value = "  -0.05%  "
$ws.Range("D19").Value = "'241.61"
$ws.Range("E19").Value = "  -3.37%  "
$ws.Range("D20").Value = "'2.082.75"
$ws.Range("E20").Value = "  -0.81%  "
$ws.Range("E21").Value = "  -1.13%  "
$ws.Range("D22").Value = "'0.9989"
$ws.Range("E22").Value = "  -0.06%  "
$ws.Range("D23").Value = "'7.519"
$ws.Range("E23").Value = "  -1.02%  "
$ws.Range("D24").Value = "'0.9997"
$ws.Range("E24").Value = "  -0.04%  "
$ws.Range("E25").Value = "  -1.19%  "
$ws.Range("D26").Value = "'159.23"
$ws.Range("E26").Value = "  -0.81%  "
$ws.Range("D27").Value = "'8.818"
$ws.Range("E27").Value = "  -0.81%  "
$ws.Range("E28").Value = "  -1.65%  "
$ws.Range("D29").Value = "'1.543"
$ws.Range("E29").Value = "  -1.51%  "
$ws.Range("D30").Value = "'4.239"
$ws.Range("E30").Value = "  -0.91%  "
$ws.Range("D31").Value = "'4.195"
$ws.Range("D32").Value = "'1.196"
$ws.Range("E32").Value = "  -0.99%  "
$ws.Range("D33").Value = "'0.05122"
$ws.Range("E33").Value = "  -1.79%  "
$ws.Range("D34").Value = "'0.7917"
$ws.Range("E34").Value = "  +4.34%  "
$ws.Range("D35").Value = "'1.865"
$ws.Range("E35").Value = "  -0.37%  "
$ws.Range("D36").Value = "'1.149"
$ws.Range("D37").Value = "'2.690"
$ws.Range("E37").Value = "  -0.72%  "
$ws.Range("D38").Value = "'1.311.05"
$ws.Range("E38").Value = "  +6.90%  "
$ws.Range("D39").Value = "'0.01871"
$ws.Range("E40").Value = "  -0.37%  "
$ws.Range("D41").Value = "'0.9449"
$ws.Range("E41").Value = "  +5.26%  "
$ws.Range("D42").Value = "'6.073"
$ws.Range("E42").Value = "  +5.84%  "
$ws.Range("D43").Value = "'107.65"
$ws.Range("E43").Value = "  -1.76%  "
$ws.Range("D44").Value = "'0.9993"
$ws.Range("E44").Value = "  -0.03%  "
$ws.Range("D45").Value = "'9.741"
$ws.Range("E45").Value = "  +2.17%  "
$ws.Range("E46").Value = "  -0.94%  "
$ws.Range("B47").Value = "RocketPoolETH"
$ws.Range("C47").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D47").Value = "'1.984.29"
$ws.Range("E47").Value = "  -0.63%  "
$ws.Range("B48").Value = "Mantle"
$ws.Range("C48").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D48").Value = "'0.5179"
$ws.Range("E48").Value = "  -0.11%  "
$ws.Range("D49").Value = "'64.25"
$ws.Range("E49").Value = "  -2.06%  "
$ws.Range("D50").Value = "'1.766"
$ws.Range("E50").Value = "  +0.02%  "
$ws.Range("D51").Value = "'7.000"
$ws.Range("E51").Value = "  -0.51%  "
